$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.005.08'
$ws.Cells.Item(2, 5).Value = '  +0.42%  '

$ws.Cells.Item(3, 4).Value = '3.067.51'
$ws.Cells.Item(3, 5).Value = '  +0.15%  '

$ws.Cells.Item(4, 5).Value = '  +0.05%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '537.42'
$ws.Cells.Item(5, 5).Value = '  -0.75%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '137.05'
$ws.Cells.Item(6, 5).Value = '  +2.46%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.999'
$ws.Cells.Item(7, 5).Value = '  +0.00%  '

$ws.Cells.Item(8, 4).Value = '3.061.21'
$ws.Cells.Item(8, 5).Value = '  +0.20%  '

$ws.Cells.Item(9, 5).Value = '  +0.48%  '

$ws.Cells.Item(10, 5).Value = '  +0.56%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.23'
$ws.Cells.Item(11, 5).Value = '  +1.89%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.453'
$ws.Cells.Item(12, 5).Value = '  -2.06%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000222'
$ws.Cells.Item(13, 5).Value = '  +0.37%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '34.40'
$ws.Cells.Item(14, 5).Value = '  -1.17%  '

$ws.Cells.Item(15, 4).Value = '3.558.27'
$ws.Cells.Item(15, 5).Value = '  +0.08%  '

$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).Value = '62.952.85'
$ws.Cells.Item(16, 5).Value = '  +0.53%  '

$ws.Cells.Item(17, 2).Value = 'TRON'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.113'
$ws.Cells.Item(17, 5).Value = '  +1.59%  '

$ws.Cells.Item(18, 4).Value = '3.061.75'
$ws.Cells.Item(18, 5).Value = '  +0.37%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.62'
$ws.Cells.Item(19, 5).Value = '  -0.97%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '469.37'
$ws.Cells.Item(20, 5).Value = '  -2.31%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.49'
$ws.Cells.Item(21, 5).Value = '  +0.29%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.694'
$ws.Cells.Item(22, 5).Value = '  -2.25%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.01'
$ws.Cells.Item(23, 5).Value = '  -2.65%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '78.37'
$ws.Cells.Item(24, 5).Value = '  -0.32%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '12.10'
$ws.Cells.Item(25, 5).Value = '  +0.12%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 5).Value = '  -0.27%  '

$ws.Cells.Item(27, 5).Value = '  -1.47%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.87'
$ws.Cells.Item(28, 5).Value = '  -5.48%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.999'
$ws.Cells.Item(29, 5).Value = '  +0.11%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '26.07'
$ws.Cells.Item(30, 5).Value = '  -0.77%  '

$ws.Cells.Item(31, 5).Value = '  +4.41%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.87'
$ws.Cells.Item(32, 5).Value = '  -2.94%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '59.09'
$ws.Cells.Item(33, 5).Value = '  +0.82%  '

$ws.Cells.Item(34, 5).Value = '  -5.73%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.46'
$ws.Cells.Item(35, 5).Value = '  +6.81%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.95'
$ws.Cells.Item(36, 5).Value = '  -0.27%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '480.59'
$ws.Cells.Item(37, 5).Value = '  -2.60%  '

$ws.Cells.Item(38, 4).Value = '3.252.52'
$ws.Cells.Item(38, 5).Value = '  +3.79%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0397'
$ws.Cells.Item(39, 5).Value = '  +0.54%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0790'
$ws.Cells.Item(40, 5).Value = '  -0.86%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.118'
$ws.Cells.Item(41, 5).Value = '  +0.62%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '8.12'
$ws.Cells.Item(42, 5).Value = '  +0.45%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.56'
$ws.Cells.Item(43, 5).Value = '  +0.22%  '

$ws.Cells.Item(44, 2).Value = 'TheGraph'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.250'
$ws.Cells.Item(44, 5).Value = '  -1.00%  '

$ws.Cells.Item(45, 2).Value = 'USDe'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.00'
$ws.Cells.Item(45, 5).Value = '  +0.09%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '123.40'
$ws.Cells.Item(46, 5).Value = '  +4.70%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '25.15'
$ws.Cells.Item(47, 5).Value = '  +1.43%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.00'
$ws.Cells.Item(48, 5).Value = '  -1.81%  '

$ws.Cells.Item(49, 5).Value = '  +0.68%  '

$ws.Cells.Item(50, 4).Value = '0.0₃0519'
$ws.Cells.Item(50, 5).Value = '  +1.10%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.01'
$ws.Cells.Item(51, 5).Value = '  -0.62%  '
